$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 158 - this shifts old rows 158:278 down to 159:279
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with the new data record
$ws.Range("A158").Value = 4
$ws.Range("B158").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C158").Value = "Los Lagos"
$ws.Range("D158").Value = 44729
$ws.Range("E158").Value = 10
$ws.Range("F158").Value = 100112003
$ws.Range("G158").Value = "Ajo"
$ws.Range("H158").Value = "Chino"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 220
$ws.Range("K158").Value = 22000
$ws.Range("L158").Value = 23000
$ws.Range("M158").Value = 22455
$ws.Range("N158").Value = "$/caja 10 kilos"
$ws.Range("O158").Value = "China"
$ws.Range("P158").Value = 2246
$ws.Range("Q158").Value = 10
$ws.Range("R158").Value = "Hortaliza"
